# This script rewrites the player roster table on Sheet1 (A1:C19) to match
# the updated data set: a reordered/refreshed list of players, their
# positions, and their teams, plus one new row (Tyler Herro moved to the end).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2..19 (A: player, B: position, C: team)
$data = @(
    @("De'Aaron Fox", "PG", "Sacramento Kings"),
    @("Dillon Brooks", "SG,SF", "Houston Rockets"),
    @("Ja Morant", "PG", "Memphis Grizzlies"),
    @("DeMar DeRozan", "SF,PF", "Sacramento Kings"),
    @("Miles Bridges", "SF,PF", "Charlotte Hornets"),
    @("Kevon Looney", "PF,C", "Golden State Warriors"),
    @("Bobby Portis", "PF,C", "Milwaukee Bucks"),
    @("Nick Richards", "C", "Phoenix Suns"),
    @("Brook Lopez", "C", "Milwaukee Bucks"),
    @("Clint Capela", "C", "Atlanta Hawks"),
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls"),
    @("Shaedon Sharpe", "SG,SF", "Portland Trail Blazers"),
    @("Scottie Barnes", "PG,SG,SF,PF", "Toronto Raptors"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls"),
    @("Luka Doncic", "PG,SG", "Dallas Mavericks"),
    @("Evan Mobley", "PG,SG,SF,PF", "Cleveland Cavaliers"),
    @("Tyler Herro", "PG,SG", "Miami Heat")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
